$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 20 (CLAVIMOX 1 GM 12 F.C.TABS.) - ratio changed from "2:0" to "1:2",
#    with price and count doubled accordingly.
$ws.Range("H20").Value = "1:2"
$ws.Range("L20").Value = 86.66
$ws.Range("N20").Value = 0.66

# 2) Insert a brand new product row before the old totals row (old row 94),
#    pushing the totals row and the footer row down by one.
$ws.Rows.Item(94).Insert()

# Copy the formatting of the row above (row 93) onto the freshly inserted
# row so the new row 94 keeps the same visual style as the rest of the table.
$ws.Range("A93:N93").Copy()
$ws.Range("A94:N94").PasteSpecial(-4122)

# Re-create the merged cells for the new data row (matches the other rows'
# B:G / H:K / L:M merge pattern).
$ws.Range("B94:G94").Merge()
$ws.Range("H94:K94").Merge()
$ws.Range("L94:M94").Merge()

# Fill in the new row's data.
$ws.Range("A94").Value = 91
$ws.Range("B94").Value = "مناديل سولو سحب"
$ws.Range("H94").Value = "30:0"
$ws.Range("L94").Value = 45
$ws.Range("N94").Value = 1

# 3) Update the running total (now on row 95) to include the new row and
#    the increased price on row 20.
$ws.Range("K95").Value = 6859.85
